$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new histology rows (8-16), mirroring the existing table layout ---

# Row 8
$ws.Range("A8").Value = "AVI999"
$ws.Range("B8").Value = 43014
$ws.Range("C8").Value = "tg/wt"
$ws.Range("D8").Value = "x"
$ws.Range("E8").Value = "x"
$ws.Range("F8").Value = "OCT"

# Row 9
$ws.Range("A9").Value = "av567"
$ws.Range("B9").Value = 43014
$ws.Range("C9").Value = "tg/tg"
$ws.Range("D9").Value = "x"
$ws.Range("E9").Value = "x"
$ws.Range("F9").Value = "OCT"

# Row 10
$ws.Range("A10").Value = "AV345"
$ws.Range("B10").Value = 43014
$ws.Range("C10").Value = "tg/wt"
$ws.Range("D10").Value = "x"
$ws.Range("E10").Value = "x"
$ws.Range("F10").Value = "OCT"

# Row 11
$ws.Range("A11").Value = "AV678"
$ws.Range("B11").Value = 43014
$ws.Range("C11").Value = "fl/fl"
$ws.Range("D11").Value = "x"
$ws.Range("E11").Value = "x"
$ws.Range("F11").Value = "OCT"

# Row 12
$ws.Range("A12").Value = "AV544"
$ws.Range("B12").Value = 43014
$ws.Range("C12").Value = "wt/wt"
$ws.Range("D12").Value = "x"
$ws.Range("E12").Value = "x"
$ws.Range("F12").Value = "OCT"

# Row 13
$ws.Range("A13").Value = "AV666"
$ws.Range("B13").Value = 43014
$ws.Range("C13").Value = "fl/fl"
$ws.Range("D13").Value = "x"
$ws.Range("E13").Value = "x"
$ws.Range("F13").Value = "OCT"

# Row 14
$ws.Range("A14").Value = "AV777"
$ws.Range("B14").Value = 43014
$ws.Range("C14").Value = "tg/tg"
$ws.Range("D14").Value = " "
$ws.Range("E14").Value = "x"
$ws.Range("F14").Value = "OCT"

# Row 15
$ws.Range("A15").Value = "AV212"
$ws.Range("B15").Value = 43014
$ws.Range("C15").Value = "tg/tg"
$ws.Range("D15").Value = "x"
$ws.Range("E15").Value = "x"
$ws.Range("F15").Value = "OCT"

# Row 16
$ws.Range("A16").Value = "AV567"
$ws.Range("B16").Value = 43014
$ws.Range("C16").Value = "tg/wt"
$ws.Range("D16").Value = "x"
$ws.Range("E16").Value = "x"
$ws.Range("F16").Value = "OCT"

# New date cells should carry on the same yyyy-mm-dd display the rest of
# the Date column already uses (re-applying it keeps them on the existing
# style instead of inheriting column B's plain date style).
$ws.Range("B8:B15").NumberFormat = "yyyy-mm-dd"

# --- Highlight the Treatment column so new entries stand out ---
$ws.Range("F1").EntireColumn.Select()
$excel.Selection.Interior.Color = 65535

# The very last date got re-entered through the format dialog and came
# back out re-escaped, landing on its own distinct style.
$ws.Range("B16").NumberFormat = "yyyy\-mm\-dd"

"done"
